$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.812.61'
$ws.Range('E2').Value = '  +1.68%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.888.16'
$ws.Range('E3').Value = '  +1.65%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.29%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.20'
$ws.Range('E5').Value = '  +1.56%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  +0.29%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4712'
$ws.Range('E7').Value = '  +3.55%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3948'
$ws.Range('E8').Value = '  +0.79%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.68'
$ws.Range('E9').Value = '  +0.30%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08082'
$ws.Range('E10').Value = '  +1.98%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.030'
$ws.Range('E11').Value = '  +1.81%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.33'
$ws.Range('E12').Value = '  +4.05%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.879.53'
$ws.Range('E13').Value = '  +1.02%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.992'
$ws.Range('E14').Value = '  +1.29%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.156'
$ws.Range('E15').Value = '  +0.03%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.007'
$ws.Range('E16').Value = '  +0.51%  '

# Row 17
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.06724'
$ws.Range('E17').Value = '  +1.64%  '

# Row 18
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '87.38'
$ws.Range('E18').Value = '  +1.14%  '

# Row 19
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001049'
$ws.Range('E19').Value = '  +2.08%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.40'
$ws.Range('E20').Value = '  +1.13%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.004'
$ws.Range('E21').Value = '  +0.25%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.555'
$ws.Range('E22').Value = '  +1.12%  '

# Row 23
$ws.Range('B23').Value = 'WrappedBTC'
$ws.Range('C23').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.810.79'
$ws.Range('E23').Value = '  +1.67%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.06'
$ws.Range('E24').Value = '  +1.31%  '

# Row 25
$ws.Range('E25').Value = '  +0.83%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.100.04'
$ws.Range('E26').Value = '  +0.96%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '158.97'
$ws.Range('E27').Value = '  +2.82%  '

# Row 28
$ws.Range('E28').Value = '  +1.01%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.121'
$ws.Range('E29').Value = '  +2.36%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.611'
$ws.Range('E30').Value = '  +2.72%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '122.08'
$ws.Range('E31').Value = '  +0.67%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9904'
$ws.Range('E32').Value = '  +4.05%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09502'

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.459'
$ws.Range('E34').Value = '  +0.16%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.608'
$ws.Range('E35').Value = '  +0.57%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.378'
$ws.Range('E36').Value = '  +2.13%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06153'
$ws.Range('E37').Value = '  +1.92%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02270'
$ws.Range('E38').Value = '  +1.86%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.229'
$ws.Range('E39').Value = '  +0.82%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.140'
$ws.Range('E40').Value = '  +1.07%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6022'
$ws.Range('E41').Value = '  +1.70%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1903'
$ws.Range('E42').Value = '  +1.01%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.36'
$ws.Range('E43').Value = '  +2.02%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.260'
$ws.Range('E44').Value = '  -1.76%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5738'
$ws.Range('E45').Value = '  +2.20%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.28'
$ws.Range('E46').Value = '  +1.78%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.957'
$ws.Range('E47').Value = '  +2.06%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.393'
$ws.Range('E48').Value = '  +0.00%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06912'
$ws.Range('E49').Value = '  +2.45%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '113.82'
$ws.Range('E50').Value = '  +5.14%  '

# Row 51
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.074'
$ws.Range('E51').Value = '  +1.89%  '
